$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 (B and C columns)
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 4

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3

$ws.Range("C5").Value = 6

$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 6

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 2

# Add new rows 8-16 with questions Q7-Q15 and their values
$ws.Range("A8").Value = "Q7"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 3

$ws.Range("A9").Value = "Q8"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 4

$ws.Range("A10").Value = "Q9"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 5

$ws.Range("A11").Value = "Q10"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 1

$ws.Range("A12").Value = "Q11"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 2

$ws.Range("A13").Value = "Q12"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 1

$ws.Range("A14").Value = "Q13"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 7

$ws.Range("A15").Value = "Q14"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 2

$ws.Range("A16").Value = "Q15"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1

# Update selection to match new extent
$ws.Range("C17").Select() | Out-Null

Write-Output "done"
